# Stage 1: update companies data
# Swap the A, B, H, I, J, K values between paired rows:
#   row 2 <-> row 3
#   row 4 <-> row 7
#   row 5 <-> row 6
# Columns C, D, E, F, G (dates/status/source/time) stay attached to their row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "H", "I", "J", "K")
$pairs = @(@(2, 3), @(4, 7), @(5, 6))

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
